# Apply "Fixed citation graph with full text" edits
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Kapil%Goyal%kapilgoyalpgi@gmail.com%0, Poonam%Chauhan%chauhan.poonamk@gmail.com%1, Komal%Chhikara%dimpichhikara@gmail.com%1, Parakriti%Gupta%parakritii@gmail.com%1, Mini P.%Singh%minipsingh@gmail.com%1]'

$ws.Range("I2").Value = ''

$ws.Range("J2").Value = 'Elsevier B.V.'

$ws.Range("E3").Value = '[Mohammed A.%Mamun%NULL%0, Mark D.%Griffiths%NULL%7]'

$ws.Range("I3").Value = ''

$ws.Range("J3").Value = 'Elsevier B.V.'

$ws.Range("E4").Value = '[Swapnajeet%Sahoo%NULL%0, Seema%Rani%NULL%1, Shaheena%Parveen%NULL%1, Ajay%Pal Singh%NULL%1, Aseem%Mehra%NULL%3, Subho%Chakrabarti%NULL%1, Sandeep%Grover%NULL%2, Cheering%Tandup%NULL%2, Cheering%Tandup%NULL%0]'

$ws.Range("I4").Value = ''

$ws.Range("J4").Value = 'Elsevier B.V.'

$ws.Range("C5").Value = 'Unknown Title'

$ws.Range("E5").Value = '[]'

$ws.Range("F5").Value = 'not found'

$ws.Range("G5").Value = 'N/A'

$ws.Range("I5").Value = ''

$ws.Range("C6").Value = 'Unknown Title'

$ws.Range("D6").Value = 'Unknown Abstract'

$ws.Range("E6").Value = '[]'

$ws.Range("F6").Value = 'not found'

$ws.Range("G6").Value = 'N/A'

$ws.Range("I6").Value = ''

$ws.Range("E7").Value = '[Chien-Cheng%Huang%NULL%0, David Hung-Tsang%Yen%hjyen@vghtpe.gov.tw%1, Hsien-Hao%Huang%NULL%1, Wei-Fong%Kao%NULL%1, Lee-Min%Wang%NULL%0, Chun-I%Huang%NULL%1, Chen-Hsen%Lee%NULL%0]'

$ws.Range("I7").Value = ''

$ws.Range("J7").Value = 'Elsevier. Published by Elsevier B.V.'

$ws.Range("C8").Value = 'Unknown Title'

$ws.Range("E8").Value = '[]'

$ws.Range("F8").Value = 'not found'

$ws.Range("G8").Value = 'N/A'

$ws.Range("H8").Value = '1970-01-01'

$ws.Range("J8").Value = ''

$ws.Range("E9").Value = '[Olaoluwa%Okusaga%NULL%0, Robert H.%Yolken%NULL%1, Patricia%Langenberg%NULL%1, Manana%Lapidus%NULL%1, Timothy A.%Arling%NULL%1, Faith B.%Dickerson%NULL%1, Debra A.%Scrandis%NULL%1, Emily%Severance%NULL%1, Johanna A.%Cabassa%NULL%1, Theodora%Balis%NULL%1, Teodor T.%Postolache%NULL%1]'

$ws.Range("I9").Value = ''

$ws.Range("J9").Value = 'Elsevier B.V.'

$val_D10 = @'
Coronavirus disease 2019 (COVID-19) was recently declared a pandemic by the WHO.
 This outbreak threatens not only physical health but also has significant repercussions on mental health.
 In recent world history, major infectious outbreaks were associated with severe mental health sequelae, including suicide.
 In this study, we systematically review the literature on suicidal outcomes during major international respiratory outbreaks, including COVID-19. We reviewed descriptive and analytic articles addressing suicide during major international respiratory outbreaks.
 We searched PubMed, Medline, Embase, Scopus, and PsycInfo databases and then utilized an independent method for study selection by a pair of reviewers.
 Two reviewers completed data abstraction and conducted a narrative summary of the findings.
 Our search generated 2,153 articles.
 Nine studies (three descriptive, five analytical, and one with mixed methodology) were eligible.
 The included studies were heterogeneous, divergent in methods, and with a low degree of evidence.
 Deducing an association between pandemics, suicide, and suicide-related outcomes remains thus poorly supported.
 Future research with better methodological characteristics, the use of longitudinal studies, and a focus on suicide as the primary outcome would allow for an in-depth understanding and formulation of the scope of this problem.

'@
$ws.Range("D10").Value = $val_D10

$ws.Range("E10").Value = '[Karine%Kahil%NULL%0, Mohamad Ali%Cheaito%NULL%1, Rawad%El Hayek%NULL%1, Marwa%Nofal%NULL%1, Sarah%El Halabi%NULL%1, Kundadak Ganesh%Kudva%NULL%1, Victor%Pereira-Sanchez%NULL%1, Samer%El Hayek%NULL%1]'

$ws.Range("I10").Value = ''

$ws.Range("J10").Value = 'Elsevier B.V.'
